# Auto-generated from the supplied OOXML diff (cryptos.xlsx "Updated
# symbol list" commit). Rewrites the Coin / Link / Price / Volume(1h)
# cells on Sheet1 that moved between the two revisions.
#
# NOTE: the Price column (D) holds numeric-looking values stored as
# plain TEXT in the workbook (t="inlineStr"), not as numbers. A bare
# Range.Value assignment of a numeric-looking string gets auto-
# converted to a real number by Excel, so for column D we prefix the
# literal with a single quote - Excel's standard "treat as text"
# quote-prefix - to keep it text, matching the source file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'244.62"
$ws.Range("D3").Value = "'21.86"
$ws.Range("D4").Value = "'5.387"
$ws.Range("D5").Value = "'0.06013"
$ws.Range("D6").Value = "'3.390"
$ws.Range("D7").Value = "'0.8153"
$ws.Range("D8").Value = "'0.9418"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "'0.1435"
$ws.Range("E9").Value = "8WazirXWRX"
$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D10").Value = "'0.07341"
$ws.Range("E10").Value = "9MandalaExchangeTokenMDX"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "'0.03421"
$ws.Range("E11").Value = "10LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03052"
$ws.Range("E12").Value = "11BitrueCoinBTR"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09413"
$ws.Range("E13").Value = "12BitMartTokenBMX"
$ws.Range("B14").Value = "MCDex"
$ws.Range("C14").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D14").Value = "'4.005"
$ws.Range("E14").Value = "13MCDexMCB"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001590"
$ws.Range("E15").Value = "14BitForexTokenBF"
$ws.Range("B16").Value = "CoinExToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D16").Value = "'0.04810"
$ws.Range("E16").Value = "15CoinExTokenCET"
$ws.Range("B17").Value = "One"
$ws.Range("C17").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D17").Value = "'0.0005942"
$ws.Range("E17").Value = "16OneONE"
$ws.Range("D18").Value = "'0.005573"
$ws.Range("D20").Value = "'0.0009893"
$ws.Range("D21").Value = "'3.669"
$ws.Range("D22").Value = "'6.431"
$ws.Range("D23").Value = "'2.172"
$ws.Range("D25").Value = "'0.1338"
$ws.Range("D40").Value = "'0.04013"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1075"
$ws.Range("E41").Value = "40BKEXTokenBKK"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.002721"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "'0.003036"
$ws.Range("E43").Value = "42KickTokenKICK"
$ws.Range("D44").Value = "'0.006640"
$ws.Range("E44").Value = "43LocalTradersLCTBestin24h"
$ws.Range("D45").Value = "'0.00005235"
$ws.Range("D48").Value = "'0.002603"
